$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ERWD-energyintensity")

# Add new rows for green hydrogen and low carbon hydrogen, mirroring the
# existing "hydrogen" row (row 11): label in column A, zeros in B:AF.
$ws.Range("A12").Value = "green hydrogen"
$ws.Range("B12:AF12").Value = 0
$ws.Range("B12:AF12").NumberFormat = "0"

$ws.Range("A13").Value = "low carbon hydrogen"
$ws.Range("B13:AF13").Value = 0
$ws.Range("B13:AF13").NumberFormat = "0"

# Make this sheet the active one/selected tab, with the two new rows selected.
$ws.Activate()
$ws.Rows("12:13").Select()
